$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1001"
$ws.Range("B2").Value = "Abdullah"
$ws.Range("C2").Value = 3452820243
$ws.Range("D2").Value = "Karachi"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2025-02-23"
$ws.Range("F2").Value = "Male"
$ws.Range("G2").Value = "Saboor"
$ws.Range("H2").Value = 3452949573
$ws.Range("I2").Value = "Gold"
$ws.Range("J2").Value = "31-03-2025"
$ws.Range("L2").Value = "Active"
$ws.Range("M2").Value = "Paid"
$ws.Range("N2").Value = "Strong"
$ws.Range("O2").Value = 70
$ws.Range("P2").Value = 6
